$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(8,9,10,12,14,15,17,18,23,34,35,36,38,40,41,43,44,49,60,61,62,64,66,67,69,70,75,86,87,88,90,92,93,95,96,101,112,113,114,116,118,119,121,122,127,138,139,140,142,144,145,147,148,153,164,167,170,174,191,194,197,201,218,221,224,228,245,248,251,255,272,275,278,282,299,302,305,309)

foreach ($r in $rows) {
    $ws.Range("G$r").Value = "System, dnasr281@gmail.com"
}
